$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.092.91"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.564.24"
$ws.Range("E3").Value = "  +0.32%  "

$ws.Range("E4").Value = "  +0.77%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.77"
$ws.Range("E5").Value = "  +1.79%  "

$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("E7").Value = "  +0.60%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.93"
$ws.Range("E8").Value = "  -0.70%  "

$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("E10").Value = "  +0.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0864"
$ws.Range("E11").Value = "  +0.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.787.23"
$ws.Range("E12").Value = "  +0.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.539.49"
$ws.Range("E13").Value = "  +1.35%  "

$ws.Range("E14").Value = "  +0.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.517"
$ws.Range("E15").Value = "  -0.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.074.01"
$ws.Range("E16").Value = "  +0.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.07"
$ws.Range("E17").Value = "  +0.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0700"
$ws.Range("E18").Value = "  -1.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.23"
$ws.Range("E19").Value = "  -1.25%  "

$ws.Range("E20").Value = "  +0.35%  "

$ws.Range("E21").Value = "  +0.49%  "

$ws.Range("E22").Value = "  +0.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.18"
$ws.Range("E23").Value = "  -0.65%  "

$ws.Range("E24").Value = "  +1.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.98"
$ws.Range("E25").Value = "  +0.48%  "

$ws.Range("E26").Value = "  -0.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.01"
$ws.Range("E27").Value = "  -0.25%  "

$ws.Range("E28").Value = "  +1.38%  "

$ws.Range("E29").Value = "  +0.58%  "

$ws.Range("E30").Value = "  +4.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0470"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("E32").Value = "  -0.09%  "

$ws.Range("E33").Value = "  +2.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.443.86"
$ws.Range("E34").Value = "  +1.07%  "

$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("E36").Value = "  -0.48%  "

$ws.Range("E37").Value = "  +1.90%  "

$ws.Range("E38").Value = "  +0.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.531"
$ws.Range("E39").Value = "  +0.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.80"
$ws.Range("E40").Value = "  +2.40%  "

$ws.Range("E41").Value = "  -0.54%  "

$ws.Range("E42").Value = "  +0.55%  "

$ws.Range("E43").Value = "  +1.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.45"
$ws.Range("E45").Value = "  -0.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.73"
$ws.Range("E46").Value = "  -0.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.702.39"
$ws.Range("E47").Value = "  +0.59%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.71"
$ws.Range("E48").Value = "  -2.27%  "

$ws.Range("E49").Value = "  +3.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0519"
$ws.Range("E50").Value = "  -0.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0956"
$ws.Range("E51").Value = "  -0.50%  "
